# "Traffic Models" slide (slide 19): expand the bullet list under
# "Each traffic can be characterized according to" with parenthetical
# clarifications, and tweak the "Consequence" paragraph's wording.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(19)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# --- Paragraph: "Packet size " -> "Packet size (bytes per transmission)" ---
$paraPacketSize = $tr.Paragraphs(6)
[void]$paraPacketSize.InsertAfter("(bytes per transmission)")

# --- Paragraph: "Interarrival time and packet queue saturation"
#     -> "Interarrival time (between successful transmissions)" ---
$paraInterarrival = $tr.Paragraphs(7)
$runTail = $paraInterarrival.Runs(2)
$runTail.Text = " "
[void]$paraInterarrival.InsertAfter("time ")
[void]$paraInterarrival.InsertAfter("(between successful transmissions)")

# --- Paragraph: "Type of packet"
#     -> "Packet queue saturation (buffer saturation level)" ---
$paraQueue = $tr.Paragraphs(8)
$runQueue1 = $paraQueue.Runs(1)
$runQueue1.Text = "Packet "
[void]$paraQueue.InsertAfter("queue ")
[void]$paraQueue.InsertAfter("saturation (buffer saturation level)")

# --- Paragraph: "Consequence: our model needs to parameterize ..."
#     -> "Consequence: Our model needs to parameterize ..." (split runs) ---
$paraConsequence = $tr.Paragraphs(10)

# "our " (chars 14-17 of the paragraph) -> "Our "
$ourRange = $paraConsequence.Characters(14, 4)
$ourRange.Text = "Our "

# Split "Consequence: " into "Consequence" + ": "
$colonRange = $paraConsequence.Characters(12, 2)
$colonRange.Text = ": "
